$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) ReadMe sheet: append the new note paragraph about the CAZyme metagenomic
#    dataset (rows 18-23; row 17 stays blank, matching the original file's
#    blank-line-between-notes convention).
# ---------------------------------------------------------------------------
$readme = $wb.Worksheets.Item("ReadMe")
$readme.Range("A18").Value = "Ashish had also re-worked his metagenomic dataset to sum up the"
$readme.Range("A19").Value = "number of CAZyme genes dedicated to a specific substrate per million"
$readme.Range("A20").Value = "reads, as opposed to last time in which he calculated a proportion"
$readme.Range("A21").Value = "of CAZyme genes dedicated to a specific substrate out of all"
$readme.Range("A22").Value = "possible CAZyme genes. I analyzed this dataset using linear mixed"
$readme.Range("A23").Value = "effect models as well."

# ---------------------------------------------------------------------------
# 2) Add a new worksheet "CAZyme metagenomic" at the end of the workbook
#    (after "litterChem") for the new dataset's LME results.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "CAZyme metagenomic"

# Re-use the header formatting (bold, centered, bordered) from the Vmax tab.
$vmax = $wb.Worksheets.Item("Vmax")
$vmax.Range("A1:E1").Copy()
$newSheet.Range("A1:E1").PasteSpecial(-4122)  # xlPasteFormats

$newSheet.Range("A1").Value = "dependent"
$newSheet.Range("B1").Value = "transformation"
$newSheet.Range("C1").Value = "Vegetation"
$newSheet.Range("D1").Value = "Precip"
$newSheet.Range("E1").Value = "interaction"

$newSheet.Range("A2").Value = "cellulose"
$newSheet.Range("C2").Value = 0.1487382805060204
$newSheet.Range("D2").Value = 0.1136010621013912
$newSheet.Range("E2").Value = "***"

$newSheet.Range("A3").Value = "chitin"
$newSheet.Range("C3").Value = 0.8149801989992265

$newSheet.Range("A4").Value = "hemicellulose"
$newSheet.Range("C4").Value = 1.657287342449503
$newSheet.Range("D4").Value = 0.310960681005965

$newSheet.Range("A5").Value = "lignin"
$newSheet.Range("C5").Value = 1.055639854891317
$newSheet.Range("D5").Value = 0.7646898136555431

$newSheet.Range("A6").Value = "oligosaccharides"
$newSheet.Range("B6").Value = "reciprocal"
$newSheet.Range("C6").Value = 1.133114552810442

$newSheet.Range("A7").Value = "peptidoglycan"

$newSheet.Range("A8").Value = "polysaccharides"

$newSheet.Range("A9").Value = "starch"
$newSheet.Range("C9").Value = 1.333409736892981

# Leave the ReadMe tab as the active/selected sheet, matching the original
# workbook's tabSelected state.
$readme.Activate()
